$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: the "NA" note no longer applies -- clear C8 back to an empty
# text cell (matching the blank C2:C7 cells above it).
$ws.Range("C8").Value = "'"
$ws.Range("C8").Style = "Normal"

# New results row appended by the script run.
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "2025-02-21"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C9").Value = "NA"
$ws.Range("D9").Value = 552
